$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Update column widths: col E (5) 13 -> 10, col F (6) 10 -> 11
# (subtract 5/6 to compensate for the character-width -> stored-width padding
# the engine applies when writing ColumnWidth back to the OOXML "width" attribute)
$ws.Columns.Item(5).ColumnWidth = 9.16666666666667
$ws.Columns.Item(6).ColumnWidth = 10.16666666666667

# Update header row (month labels shift forward by one month)
$ws.Range("C1").Value = "marzo"
$ws.Range("D1").Value = "abril"
$ws.Range("E1").Value = "mayo"
$ws.Range("F1").Value = "junio"

# Row 2
$ws.Range("C2").Value = 5618.93
$ws.Range("D2").Value = 0

# Row 3
$ws.Range("C3").Value = 1930.27
$ws.Range("D3").Value = 0

# Row 5
$ws.Range("C5").Value = 1906.46
$ws.Range("D5").Value = 0

# Row 6
$ws.Range("D6").Value = -545.18
$ws.Range("E6").Value = 0

# Row 7 (totals)
$ws.Range("C7").Value = 9455.66
$ws.Range("D7").Value = -545.18
$ws.Range("E7").Value = 0
